# Sửa lại template báo cáo tổng hợp cơ sở
# Xoá dòng "Phụ cấp tại LONG XUYÊN" (row 13) và cập nhật lại các giá trị
# tính toán do thay đổi chiến lược chạy multi process.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Lương")

# Xoá toàn bộ dòng 13 ("Phụ cấp tại LONG XUYÊN") - các dòng bên dưới tự dịch lên 1
$ws.Rows("13:13").Delete()

# Cập nhật các giá trị bị thay đổi sau khi tính lại lương
$ws.Range("B12").Value = 3
$ws.Range("B13").Value = 535714.2857142857
$ws.Range("B32").Value = 535714.2857142857
$ws.Range("B34").Value = 535714.2857142857
